$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

# Updated LTP (col B) and PREV (col C) prices for 6 Jun 24 (algo shares fix)
$data = @(
    @(2, 587, 582.95),
    @(3, 8268.1, 8055),
    @(4, 3186.1, 3135.85),
    @(5, 550.75, 547.85),
    @(6, 238.13, 235.4),
    @(7, 1975.25, 1973.4),
    @(8, 7703, 7703),
    @(9, 195.28, 198.74),
    @(10, 280.45, 279.7),
    @(11, 248.91, 247.8),
    @(12, 53381.2, 53435.15),
    @(13, 14195.8, 13802.95),
    @(14, 913.75, 895.15),
    @(15, 4972.65, 5026.25),
    @(16, 4209.2, 4268.8),
    @(17, 197.12, 196.73),
    @(18, 1816.5, 1796.1),
    @(19, 761.55, 756.2),
    @(20, 558.4, 558.55),
    @(21, 1409.7, 1447.6),
    @(22, 1035.35, 1039.65),
    @(23, 659.25, 662.45),
    @(24, 3165.5, 3094.9),
    @(25, 330.1, 338.35),
    @(26, 25969.25, 25990.15),
    @(27, 440.1, 443.2),
    @(28, 291.95, 297.6),
    @(29, 556.8, 554.5),
    @(30, 796.95, 787.9),
    @(31, 849.65, 833.1),
    @(32, 965.2, 974.65),
    @(33, 481.25, 482.6),
    @(34, 167.03, 168.55),
    @(35, 516.15, 512.65)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Move the active selection to match the saved cursor position
$ws.Range("N22").Select()
